# NIT-9005464445.xlsx — "Actualiza base de datos EC y agrega parte 1 de
# nuevos estado de cuenta"
#
# The workbook previously listed two workers under the NIT:
#   row16: CC 1047450387 - LUIS MANUEL MERCADO BORNACHERA - periodo 2507
#   row17: CC 73095854   - DANIEL ALBERTO LYNTON ELLES     - periodo 2207
#
# The update drops the first worker (LUIS MANUEL...) from this
# "estado de cuenta" and keeps only DANIEL ALBERTO's record, now as the
# single data row. The header counters (Cant. Trabajadores / Cant.
# Periodos) and the "Valor Mora" summary field are updated accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Valor Mora summary (top of the sheet) now mirrors the remaining
# worker's "Periodo Mora" value.
$ws.Range("E11").Value = 33333

# Cant. Trabajadores / Cant. Periodos go from 2 to 1.
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1

# Overwrite row 16 (currently LUIS MANUEL's record) with what used to be
# row 17's data (DANIEL ALBERTO's record) so it keeps row 16's existing
# formatting/styles.
$ws.Range("C16").Value = "73095854"
$ws.Range("D16").Value = "DANIEL ALBERTO LYNTON ELLES"
$ws.Range("E16").Value = "2207"
$ws.Range("F16").Value = 33333
$ws.Range("G16").Value = 1000000

# Remove the now-duplicated old row 17, shifting the signature block
# (old rows 22/23) up to rows 21/22.
$ws.Rows("17:17").Delete()
